$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New game results to append (Away team, Away Pts, Home team, Home Pts, Overtime, Attend., Arena, Win, Loss)
$games = @(
    @("Philadelphia 76ers", 122, "Indiana Pacers", 134, "No", 17832, "Gainbridge Fieldhouse", "Indiana Pacers", "Philadelphia 76ers"),
    @("Utah Jazz", 123, "Washington Wizards", 108, "No", 17832, "Capital One Arena", "Utah Jazz", "Washington Wizards"),
    @("Minnesota Timberwolves", 96, "Brooklyn Nets", 94, "No", 17832, "Barclays Center", "Minnesota Timberwolves", "Brooklyn Nets"),
    @("Boston Celtics", 143, "Miami Heat", 110, "No", 17832, "Kaseya Center", "Boston Celtics", "Miami Heat"),
    @("Denver Nuggets", 84, "New York Knicks", 122, "No", 17832, "Madison Square Garden (IV)", "New York Knicks", "Denver Nuggets"),
    @("Sacramento Kings", 134, "Golden State Warriors", 133, "No", 17832, "Chase Center", "Sacramento Kings", "Golden State Warriors"),
    @("Chicago Bulls", 132, "Los Angeles Lakers", 141, "No", 17832, "Crypto.com Arena", "Los Angeles Lakers", "Chicago Bulls")
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$r = $lastRow + 1

foreach ($game in $games) {
    $ws.Cells.Item($r, 1).Value = $game[0]
    $ws.Cells.Item($r, 2).Value = $game[1]
    $ws.Cells.Item($r, 3).Value = $game[2]
    $ws.Cells.Item($r, 4).Value = $game[3]
    $ws.Cells.Item($r, 5).Value = $game[4]
    $ws.Cells.Item($r, 6).Value = $game[5]
    $ws.Cells.Item($r, 7).Value = $game[6]
    $ws.Cells.Item($r, 8).Value = $game[7]
    $ws.Cells.Item($r, 9).Value = $game[8]

    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 2)).NumberFormat = "#,##0"
    $ws.Range($ws.Cells.Item($r, 4), $ws.Cells.Item($r, 4)).NumberFormat = "#,##0"

    $r = $r + 1
}
